# DaySale_2025-06-18_00-00.xlsx - "Upload new version with timestamp"
#
# Net effect of the re-upload on the sheet data: the "COLA - CHOND 30 TABS"
# row (row 8) picked up fresher numbers -
#   H8 (الرصيد الحالي / current balance)   0:0  -> 0:1
#   P8 (سعر البيع / selling price)          82.50 -> 41.25
#   Q8 (عدد التعاملات / # of transactions)  0:2  -> 0:1
# and the grand total in P13 was refreshed to match (527.67 -> 486.42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H8 / Q8 already match an existing text value ("0:1") elsewhere in the
# sheet, so a plain text assignment is fine.
$ws.Range("H8").Value = "0:1"
$ws.Range("Q8").Value = "0:1"

# P8 is a text-typed cell (numFmt "0.00" but stored as a string "82.5000"),
# so force text formatting before writing the new value to keep it a
# string instead of Excel auto-converting it to a number, then restore the
# original numeric display format.
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "41.2500"
$ws.Range("P8").NumberFormat = "0.00"

# Refresh the total shown in the merged P13:Q13 footer cell.
$ws.Range("P13").Value = 486.42
